$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in header P1: vecation_day -> vacation_day
$ws.Range("P1").Value = "vacation_day"

# Row 6: clear the entry/exit time + name/sum fields that are no longer populated,
# bump the day value, and record a vacation day instead.
$ws.Range("A6").Value = 19
$ws.Range("C6").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("F6").ClearContents()
$ws.Range("G6").ClearContents()
$ws.Range("I6").ClearContents()
$ws.Range("O6").ClearContents()
$ws.Range("P6").Value = 1
